$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that must stay TEXT (e.g. "0.4590", "27.200.71").
# Forcing NumberFormat to Text ("@") before the assignment stops Excel from
# re-parsing the string as a number (which would drop trailing zeros / reformat
# thousands separators). Resetting the Style back to "Normal" afterwards removes
# the now-unneeded text format so the cell style matches the original workbook.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.200.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.33%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.720.49"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.65%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.70%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4590"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.29%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3427"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.47%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "41.98"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07253"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.62%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.044"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.77"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.830"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.720.87"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.27%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.862"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.57%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.56"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.83%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001040"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06332"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.55%  "
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.611"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.224.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.84"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.129"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.77"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.920.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.134"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.45"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.024"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09114"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.597"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.341"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.93%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02194"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.97%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05845"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "11.04"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.1998"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.756"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.404"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.94%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5906"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.131"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.470"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.63%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.72"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.587"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5635"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "119.06"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.862"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.20%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06657"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.58%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.084"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.001"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.07%  "
